# Automatische test-sync: 2025-08-05 18:15:50
# Adds the 5th "Wil je deze klant bellen?" test mail entry to the Logs sheet
# and bumps the corresponding Dashboard count.

$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append row 26 -------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 26
$logs.Cells.Item($newRow, 1).Value2  = "Wil je deze klant bellen?"
$logs.Cells.Item($newRow, 2).Value2  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value2  = "Testmail #5: Wil je deze klant bellen?"
$logs.Cells.Item($newRow, 4).Value2  = "Klantenservice / Contact"
$logs.Cells.Item($newRow, 5).Value2  = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value2  = "2025-08-05 18:15:12"
$logs.Cells.Item($newRow, 7).Value2  = "Ja"
$logs.Cells.Item($newRow, 8).Value2  = "Ja"
$logs.Cells.Item($newRow, 9).Value2  = "Nee"
$logs.Cells.Item($newRow, 10).Value2 = "Nee"

# --- Extend conditional formatting ranges from row 25 to the new row 26 ---------
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "25")
    $newRange = $logs.Range($col + "2:" + $col + "26")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Sheet "Dashboard": bump the "Klantenservice / Contact" count --------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value2 = 4
